$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: append extra (red) text after the existing sentence.
#    "This is a Microsoft word document." (34 chars) ->
#      run1 (plain):  "This is a Microsoft word document.  "
#      run2 (red):    "(This is a change - Ve"
#      run3 (red):    "rsion for main branch"
#      run4 (red):    ")"
# ---------------------------------------------------------------------------
$pos = $d.Paragraphs(1).Range.End - 1   # position right before the 1st pilcrow

$r = $d.Range($pos, $pos)
$r.InsertAfter("  ")
$pos = $pos + 2

$seg = "(This is a change " + [char]0x2013 + " Ve"
$r = $d.Range($pos, $pos)
$r.InsertAfter($seg)
$r2 = $d.Range($pos, $pos + $seg.Length)
$r2.Font.Color = 255
$pos = $pos + $seg.Length

$seg = "rsion for main branch"
$r = $d.Range($pos, $pos)
$r.InsertAfter($seg)
$r2 = $d.Range($pos, $pos + $seg.Length)
$r2.Font.Color = 255
$pos = $pos + $seg.Length

$seg = ")"
$r = $d.Range($pos, $pos)
$r.InsertAfter($seg)
$r2 = $d.Range($pos, $pos + $seg.Length)
$r2.Font.Color = 255
$pos = $pos + $seg.Length

# ---------------------------------------------------------------------------
# 2. Remove the trailing "ank God almighty, we are free at last." paragraph
#    (the paragraph that used the NormalWeb style), keeping the
#    "Shall be lifted-nevermore!" paragraph right before it.
# ---------------------------------------------------------------------------
$found = $true
while ($found) {
    $found = $false
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text -like "*God almighty, we are free at last.*") {
            $para.Range.Delete()
            $found = $true
            break
        }
    }
}
